{"js": "// The author dropped the \"Statistics ... \u0438\u0434\u0435\u044f \u0437\u0430 \u0440\u0430\u0437\u0448\u0438\u0440\u044f\u0432\u0430\u043d\u0435 \u043d\u0430 \u043f\u0440\u043e\u0435\u043a\u0442\u0430\n// \u0432 \u0431\u044a\u0434\u0435\u0449\u0435.\" sentence (a parked idea for a future admin-console feature)\n// from the project write-up. Removing that whole paragraph (mark and\n// all) merges it away; Word then leaves its \"last edit\" bookmark\n// (_GoBack) sitting at the point where the edit happened, i.e. right at\n// the start of the paragraph that now follows.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph that carries the dropped sentence.\nlet targetIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"Statistics\") !== -1 && text.indexOf(\"\u0431\u044a\u0434\u0435\u0449\u0435\") !== -1) {\n    targetIndex = i;\n    break;\n  }\n}\n\nif (targetIndex !== -1) {\n  // Drop the pre-existing _GoBack bookmark (Word only ever keeps one); it\n  // will be re-created at the edit location below.\n  context.document.deleteBookmark(\"_GoBack\");\n\n  paragraphs.items[targetIndex].delete();\n  await context.sync();\n\n  // After the delete, the paragraph that used to follow now sits at\n  // targetIndex; that's where editing last happened.\n  const afterParagraphs = body.paragraphs;\n  afterParagraphs.load(\"items\");\n  await context.sync();\n\n  if (targetIndex < afterParagraphs.items.length) {\n    const landingParagraph = afterParagraphs.items[targetIndex];\n    landingParagraph.getRange(\"Start\").insertBookmark(\"_GoBack\");\n  }\n\n  await context.sync();\n}\n", "ps1": "# The author dropped the \"Statistics ... \u0438\u0434\u0435\u044f \u0437\u0430 \u0440\u0430\u0437\u0448\u0438\u0440\u044f\u0432\u0430\u043d\u0435 \u043d\u0430 \u043f\u0440\u043e\u0435\u043a\u0442\u0430\n# \u0432 \u0431\u044a\u0434\u0435\u0449\u0435.\" sentence (a parked idea for a future admin-console feature)\n# from the project write-up. Removing that whole paragraph (mark and\n# all) merges it away; Word then leaves its \"last edit\" bookmark\n# (_GoBack) sitting at the point where the edit happened, i.e. right at\n# the start of the paragraph that now follows.\n\n$d = $word.ActiveDocument\n\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*Statistics*\" -and $t -like \"*\u0431\u044a\u0434\u0435\u0449\u0435*\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -ne -1) {\n    # Word keeps a single _GoBack bookmark; drop the old one, it gets\n    # re-created below at the spot where the edit actually happened.\n    if ($d.Bookmarks.Exists(\"_GoBack\")) {\n        $d.Bookmarks(\"_GoBack\").Delete()\n    }\n\n    # Deleting the paragraph's Range removes its text AND its paragraph\n    # mark, merging the paragraph away entirely.\n    $d.Paragraphs.Item($targetIndex).Range.Delete()\n\n    # The paragraph that used to follow now occupies $targetIndex; that\n    # is where editing last happened, so _GoBack lands at its start.\n    if ($targetIndex -le $d.Paragraphs.Count) {\n        $landing = $d.Paragraphs.Item($targetIndex).Range\n        $landing.Collapse(1)\n        $d.Bookmarks.Add(\"_GoBack\", $landing)\n    }\n}\n"}
